# "updated class Edit with DataDriven"
#
# The "Class" sheet gains a small, hand-entered data-driven test grid
# (edit-record fields + notes/recording columns), the previously-selected
# "Program" tab loses focus to "Class", and the now-unused bespoke
# 14pt-Arial style on the old Class sheet is cleared back to Normal.

$wb = $excel.ActiveWorkbook

$wsProgram = $wb.Worksheets.Item("Program")
$wsClass   = $wb.Worksheets.Item("Class")

# --- Program sheet: selection moves, it's no longer the active tab ---
$wsProgram.Range("C7").Select()

# --- Class sheet: clear the bespoke font/style left on the header rows ---
# (row 1-2 previously carried explicit fonts; the edit drops back to the
# workbook default "Normal" style and lets row height auto-fit again)
$wsClass.Range("A1:D2").Style = "Normal"
$wsClass.Rows("1:2").AutoFit()

# --- Class sheet: new data-driven edit test rows/columns ---
$wsClass.Range("A3").Value = "ValidEditData"
$wsClass.Range("D3").Value = "Saranya M"

$wsClass.Range("E1").Value = "ClassDesc"
$wsClass.Range("A4").Value = "OptionalValid"

$wsClass.Range("F1").Value = "Comments"
$wsClass.Range("F4").Value = "Playwright with Appium"

$wsClass.Range("G1").Value = "Notes"
$wsClass.Range("G4").Value = "c:/user/Notes"

$wsClass.Range("H4").Value = "c:/Recordings"
$wsClass.Range("H1").Value = "Recording"

$wsClass.Range("A5").Value = "NumericOrAlphaData"

$wsClass.Range("E4").Value = "Playwright"
$wsClass.Range("E5").Value = 23423
$wsClass.Range("F5").Value = 435254
$wsClass.Range("G5").Value = 23452
$wsClass.Range("H5").Value = 123515

# --- Class sheet becomes the active tab/selected cell ---
$wsClass.Activate()
$wsClass.Range("D10").Select()
